# Apply "data till 14 Apr 10AM" update
# - New sales entries in Orders / Collection sheets (rows with new daily figures)
# - New cell comments left by Vijay
# - AmtToCollect visibility/filters follow automatically from the formulas

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Orders sheet: add the newly recorded per-day sale formulas (qty * 1500)
# ---------------------------------------------------------------------------
$orders = $wb.Worksheets.Item("Orders")

$orders.Range("R10").Formula  = "=6*1500"

$orders.Range("U21").Formula  = "=2*1500"

$orders.Range("U27").Formula  = "=2*1500"

$orders.Range("U34").Formula  = "=1*1500"

$orders.Range("U36").Formula  = "=5*1500"

$orders.Range("U39").Formula  = "=3*1500"

$orders.Range("U41").Formula  = "=2*1500"

$orders.Range("U42").Formula  = "=2*1500"

$orders.Range("Q43").Formula  = "=1*1500"
$orders.Range("U43").Formula  = "=3*1500"

$orders.Range("U45").Formula  = "=4*1500"

$orders.Range("V62").Formula  = "=2*1500"

$orders.Range("Q80").Formula  = "=2*1500"
$orders.Range("U80").Formula  = "=1*1500"
$orders.Range("V80").Formula  = "=1*1500"

$orders.Range("Q85").Formula  = "=1*1500"

# ---------------------------------------------------------------------------
# Collection sheet: record the matching amounts actually collected
# ---------------------------------------------------------------------------
$collection = $wb.Worksheets.Item("Collection")

$collection.Range("Q10").Value = 9000

$collection.Range("T36").Value = 4000

$collection.Range("P43").Value = 1500

$collection.Range("R58").Value = 2000

$collection.Range("U62").Value = 3000

$collection.Range("P80").Value = 3000
$collection.Range("T80").Value = 1500
$collection.Range("U80").Value = 1500

$collection.Range("P85").Value = 1500

# ---------------------------------------------------------------------------
# Comments left by Vijay documenting the new entries
# ---------------------------------------------------------------------------
$collection.Range("T36").AddComment("Vijay:`n3000-Digital`n1000-Cash")
$collection.Range("P80").AddComment("Vijay:`n2 Mobile sold to FOS")
$collection.Range("T80").AddComment("Vijay:`n1 mobile sold to FOS")
$collection.Range("U80").AddComment("Vijay:`nSold to customer for Rs 1650")
